# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", holding
#   the per-fund holdings snapshot for the new quarter (same layout as the
#   "2021-Q4" sheet).
# - Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data (pushing the existing "2021-Q4" row down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the "2022-Q1" worksheet right after "2021-Q4" (and therefore
#    right before "总计").
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(1)
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

function Set-HeaderCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1         # xlContinuous
}

function Set-IndexCell($ws, $row, $col, $number) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $number
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1         # xlContinuous
}

function Set-TextCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
}

function Set-NumberCell($ws, $row, $col, $number) {
    $ws.Cells.Item($row, $col).Value = $number
}

# Header row
Set-HeaderCell $q1 1 2 "基金代码"
Set-HeaderCell $q1 1 3 "基金名称"
Set-HeaderCell $q1 1 4 "基金规模"
Set-HeaderCell $q1 1 5 "股票总仓位"
Set-HeaderCell $q1 1 6 "仓位占比"
Set-HeaderCell $q1 1 7 "持有市值(亿元)"
Set-HeaderCell $q1 1 8 "仓位排名"

# Row 2 - 银河君尚灵活配置混合A
Set-IndexCell $q1 2 1 0
Set-TextCell  $q1 2 2 "519613"
$q1.Cells.Item(2, 3).Value = "银河君尚灵活配置混合A"
Set-TextCell  $q1 2 4 "5.70"
Set-TextCell  $q1 2 5 "29.99"
Set-TextCell  $q1 2 6 "0.49"
Set-TextCell  $q1 2 7 "0.0279"
Set-NumberCell $q1 2 8 5

# Row 3 - 银河君尚灵活配置混合I
Set-IndexCell $q1 3 1 1
Set-TextCell  $q1 3 2 "519615"
$q1.Cells.Item(3, 3).Value = "银河君尚灵活配置混合I"
Set-TextCell  $q1 3 4 "5.70"
Set-TextCell  $q1 3 5 "29.99"
Set-TextCell  $q1 3 6 "0.49"
Set-TextCell  $q1 3 7 "0.0279"
Set-NumberCell $q1 3 8 5

# Row 4 - 银河君尚灵活配置混合C
Set-IndexCell $q1 4 1 2
Set-TextCell  $q1 4 2 "519614"
$q1.Cells.Item(4, 3).Value = "银河君尚灵活配置混合C"
Set-TextCell  $q1 4 4 "0.23"
Set-TextCell  $q1 4 5 "29.99"
Set-TextCell  $q1 4 6 "0.49"
Set-TextCell  $q1 4 7 "0.0011"
Set-NumberCell $q1 4 8 5

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q1" summary row into "总计", ahead of the
#    existing "2021-Q4" row, and renumber the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
# Excel's row-insert copies formatting down from the row above (the bold
# header); strip that back to the sheet default before (re)styling the
# cells we actually want styled.
$total.Range("A2:D2").ClearFormats()

Set-IndexCell $total 2 1 0
$total.Cells.Item(2, 2).NumberFormat = "@"
$total.Cells.Item(2, 2).Value = "2022-Q1"
Set-NumberCell $total 2 3 3
Set-NumberCell $total 2 4 0.06

# Renumber the pre-existing "2021-Q4" row, now shifted to row 3.
Set-IndexCell $total 3 1 1
